$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: the "CreateNeonUser" request body now also carries USER_TRUID, and the
# row got taller to fit the extra wrapped text.
$ws.Range("I2").Value = "USER_NAME=Neon_JDRUser4@1p.com||USER_PASSWORD=1234qwer`$`$!||PASSWORD_GENERATE=NO||EMAIL_GENERATE=YES||USER_FIRST_NAME=JANARDHAN4||USER_LAST_NAME=E4||USER_MIDDLE_NAME=REDDY4||USER_TRUID=c32994ec-6dcd-4884-ab42-682bbc0f9e8b"
$ws.Range("A2:M2").RowHeight = 60

# Row 5 (UpdateUserStatus / OPQA-DDD): the test now expects the call to fail,
# so the description and the expected validation result both change.
$ws.Range("B5").Value = "Verify that based on truid, user status can't be updated and check the error status using STeAM API"
$ws.Range("K5").Value = "status=200||rc=55003||fn[1].error=Action UPDATE_NEON_USER_STATUS is not supported"

# Row 6 (GetLoginUNP): its dependency test now points at OPQA-AAA instead of
# OPQA-DDD, and the expected validation drops the now-unused isActive check.
$ws.Range("J6").Value = "OPQA-AAA"
$ws.Range("K6").Value = "status=200||rc=OK||User.userID=(OPQA-AAA_user.userID)||User.truId=(OPQA-AAA_user.truID)||UserInfo.USER_INFO_FIRST_NAME=(OPQA-CCC_UserInfo.USER_INFO_FIRST_NAME)||UserInfo.USER_INFO_MIDDLE_NAME=(OPQA-CCC_UserInfo.USER_INFO_MIDDLE_NAME)||UserInfo.USER_INFO_LAST_NAME=(OPQA-CCC_UserInfo.USER_INFO_LAST_NAME)"

# Rows 6-8: the BODY column (I) now wraps text like the rest of the sheet.
$ws.Range("I6").WrapText = $true
$ws.Range("I7").WrapText = $true
$ws.Range("I8").WrapText = $true
